$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.722.39"
$ws.Range("E2").Value = "  -4.15%  "
$ws.Range("D3").Value = "1.817.16"
$ws.Range("E3").Value = "  -3.05%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'277.79"
$ws.Range("E5").Value = "  -7.80%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.5098"
$ws.Range("E7").Value = "  -5.12%  "
$ws.Range("D8").Value = "'0.3528"
$ws.Range("E8").Value = "  -5.82%  "
$ws.Range("D9").Value = "'44.55"
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("D10").Value = "'0.06675"
$ws.Range("E10").Value = "  -7.21%  "
$ws.Range("D11").Value = "'20.06"
$ws.Range("E11").Value = "  -7.11%  "
$ws.Range("D12").Value = "'0.8267"
$ws.Range("E12").Value = "  -7.12%  "
$ws.Range("D13").Value = "'0.07903"
$ws.Range("E13").Value = "  -3.19%  "
$ws.Range("D14").Value = "1.802.97"
$ws.Range("E14").Value = "  -3.85%  "
$ws.Range("D15").Value = "'5.074"
$ws.Range("E15").Value = "  -4.64%  "
$ws.Range("D16").Value = "'87.67"
$ws.Range("E16").Value = "  -6.19%  "
$ws.Range("D17").Value = "'1.0000"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'14.08"
$ws.Range("E18").Value = "  -5.30%  "
$ws.Range("D19").Value = "'0.000008031"
$ws.Range("E19").Value = "  -5.95%  "
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "25.770.69"
$ws.Range("E21").Value = "  -4.12%  "
$ws.Range("E22").Value = "  -4.99%  "
$ws.Range("D23").Value = "'9.997"
$ws.Range("E23").Value = "  -5.96%  "
$ws.Range("D24").Value = "'6.106"
$ws.Range("E24").Value = "  -4.69%  "
$ws.Range("D25").Value = "'2.225"
$ws.Range("E25").Value = "  -2.93%  "
$ws.Range("D26").Value = "'141.66"
$ws.Range("E26").Value = "  -3.24%  "
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").Value = "'17.10"
$ws.Range("E28").Value = "  -5.50%  "
$ws.Range("D29").Value = "'109.32"
$ws.Range("E29").Value = "  -4.15%  "
$ws.Range("D30").Value = "'4.328"
$ws.Range("E30").Value = "  -8.36%  "
$ws.Range("D31").Value = "'4.236"
$ws.Range("E31").Value = "  -8.17%  "
$ws.Range("D32").Value = "'0.08770"
$ws.Range("E32").Value = "  -4.12%  "
$ws.Range("D33").Value = "'0.04890"
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("D34").Value = "'0.7278"
$ws.Range("E34").Value = "  -10.58%  "
$ws.Range("D35").Value = "'1.136"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").Value = "'2.864"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("D37").Value = "'0.9993"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").Value = "'3.137"
$ws.Range("E38").Value = "  -2.49%  "
$ws.Range("D39").Value = "'2.383"
$ws.Range("E39").Value = "  -9.46%  "
$ws.Range("D40").Value = "'0.01854"
$ws.Range("E40").Value = "  -5.20%  "
$ws.Range("D41").Value = "'0.5165"
$ws.Range("E41").Value = "  -14.06%  "
$ws.Range("D42").Value = "'0.9659"
$ws.Range("E42").Value = "  -9.85%  "
$ws.Range("D43").Value = "'6.220"
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("D44").Value = "'110.99"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("D45").Value = "'8.026"
$ws.Range("E45").Value = "  -9.79%  "
$ws.Range("D46").Value = "'1.0000"
$ws.Range("D47").Value = "'0.4574"
$ws.Range("E47").Value = "  -10.56%  "
$ws.Range("E48").Value = "  -8.89%  "
$ws.Range("D49").Value = "'36.48"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("E50").Value = "  -8.22%  "
$ws.Range("E51").Value = "  -8.65%  "
